$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 1.007722494170824
$ws.Range("C9").Value = 1.007352044600634
$ws.Range("C10").Value = 1.015105312030071
$ws.Range("C11").Value = 1.014546911008609
$ws.Range("C12").Value = 1.014426998693669
$ws.Range("C13").Value = 1.014446296999169
$ws.Range("C14").Value = 1.007324572005172
$ws.Range("C15").Value = 1.015034496969908
$ws.Range("C16").Value = 1.013406994424825
$ws.Range("C17").Value = 1.013018847194835
$ws.Range("C18").Value = 1.013542769290194
$ws.Range("C19").Value = 1.013108967955256
